$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = theta_se : replace (nan) placeholders with actual standard errors
$ws.Range("B4").Value = "(0.28)"
$ws.Range("C4").Value = "(0.17)"
$ws.Range("D4").Value = "(0.39)"
$ws.Range("E4").Value = "(0.61)"
$ws.Range("F4").Value = "(0.83)"
$ws.Range("G4").Value = "(0.12)"
$ws.Range("H4").Value = "(1.0)"
$ws.Range("I4").Value = "(0.68)"
$ws.Range("J4").Value = "(0.36)"
$ws.Range("K4").Value = "(0.18)"
$ws.Range("L4").Value = "(0.38)"

# Row 6 = lambda_se : replace (nan) placeholders with actual standard errors
$ws.Range("B6").Value = "(0.16)"
$ws.Range("C6").Value = "(0.1)"
$ws.Range("D6").Value = "(0.29)"
$ws.Range("E6").Value = "(0.14)"
$ws.Range("F6").Value = "(0.28)"
$ws.Range("G6").Value = "(0.64)"
$ws.Range("H6").Value = "(0.71)"
$ws.Range("I6").Value = "(0.07)"
$ws.Range("J6").Value = "(0.16)"
$ws.Range("K6").Value = "(0.5)"
$ws.Range("L6").Value = "(0.97)"
